$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("N6").Value = 8
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 2.75

# Row 7
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 13

# Row 11
$ws.Range("G11").Value = 1.73
$ws.Range("I11").Value = 4.2
$ws.Range("J11").Value = 2.25
$ws.Range("Q11").Value = 1.53
$ws.Range("R11").Value = 2.4
$ws.Range("AJ11").Value = 15
$ws.Range("BB11").Value = 151

# Row 13
$ws.Range("M13").Value = 1.02
$ws.Range("N13").Value = 21
$ws.Range("AC13").Value = 23
$ws.Range("AS13").Value = 81
$ws.Range("AV13").Value = 34
$ws.Range("AZ13").Value = 29
$ws.Range("BB13").Value = 81

# Row 14
$ws.Range("N14").Value = 26
$ws.Range("O14").Value = 1.08
$ws.Range("P14").Value = 8

# Row 17
$ws.Range("O17").Value = 1.25
$ws.Range("P17").Value = 3.75
$ws.Range("Q17").Value = 1.83
$ws.Range("R17").Value = 2.03
$ws.Range("S17").Value = 1.36
$ws.Range("T17").Value = 3
$ws.Range("AE17").Value = 19
$ws.Range("AK17").Value = 67
$ws.Range("AT17").Value = 3
$ws.Range("BD17").Value = 126

# Row 19
$ws.Range("J19").Value = 3
$ws.Range("K19").Value = 2.25
$ws.Range("L19").Value = 3.4
$ws.Range("M19").Value = 1.04
$ws.Range("N19").Value = 13
$ws.Range("O19").Value = 1.22
$ws.Range("P19").Value = 4
$ws.Range("Q19").Value = 1.73
$ws.Range("R19").Value = 2.08
$ws.Range("S19").Value = 1.33
$ws.Range("T19").Value = 3.25
$ws.Range("U19").Value = 1.62
$ws.Range("V19").Value = 2.2
$ws.Range("W19").Value = 10
$ws.Range("AC19").Value = 13
$ws.Range("AG19").Value = 151
$ws.Range("AH19").Value = 11
$ws.Range("AI19").Value = 15
$ws.Range("AO19").Value = 13
$ws.Range("AP19").Value = 21
$ws.Range("AS19").Value = 126
$ws.Range("AT19").Value = 3.25
$ws.Range("AU19").Value = 7.5
$ws.Range("AZ19").Value = 51
$ws.Range("BD19").Value = 126

# Row 21
$ws.Range("H21").Value = 4.33
$ws.Range("J21").Value = 5
$ws.Range("K21").Value = 2.5
$ws.Range("L21").Value = 2.1
$ws.Range("N21").Value = 19
$ws.Range("O21").Value = 1.14
$ws.Range("P21").Value = 5.5
$ws.Range("Q21").Value = 1.5
$ws.Range("R21").Value = 2.5
$ws.Range("S21").Value = 1.25
$ws.Range("T21").Value = 3.75
$ws.Range("U21").Value = 1.57
$ws.Range("V21").Value = 2.25
$ws.Range("W21").Value = 19
$ws.Range("AB21").Value = 34
$ws.Range("AC21").Value = 19
$ws.Range("AD21").Value = 8.5
$ws.Range("AE21").Value = 15
$ws.Range("AF21").Value = 41
$ws.Range("AG21").Value = 126
$ws.Range("AH21").Value = 10
$ws.Range("AI21").Value = 9.5
$ws.Range("AL21").Value = 12
$ws.Range("AP21").Value = 26
$ws.Range("AQ21").Value = 81
$ws.Range("AR21").Value = 81
$ws.Range("AS21").Value = 151
$ws.Range("AT21").Value = 3.75
$ws.Range("BA21").Value = 41
$ws.Range("BC21").Value = 351

# Row 22
$ws.Range("G22").Value = 1.3
$ws.Range("H22").Value = 5.75
$ws.Range("I22").Value = 7.5
$ws.Range("J22").Value = 1.73
$ws.Range("K22").Value = 2.88
$ws.Range("L22").Value = 6.5
$ws.Range("U22").Value = 1.57
$ws.Range("V22").Value = 2.25
$ws.Range("X22").Value = 9.5
$ws.Range("Y22").Value = 9
$ws.Range("AD22").Value = 12
$ws.Range("AE22").Value = 17
$ws.Range("AG22").Value = 126
$ws.Range("AJ22").Value = 21
$ws.Range("AL22").Value = 41
$ws.Range("AQ22").Value = 15
$ws.Range("AZ22").Value = 101

# Row 23
$ws.Range("M23").Value = 1.07
$ws.Range("N23").Value = 9
$ws.Range("O23").Value = 1.36
$ws.Range("P23").Value = 3
$ws.Range("Q23").Value = 2.2
$ws.Range("R23").Value = 1.65

# Row 29
$ws.Range("G29").Value = 3.7
$ws.Range("H29").Value = 3.45
$ws.Range("I29").Value = 1.93
$ws.Range("J29").Value = 4.2
$ws.Range("K29").Value = 2.12
$ws.Range("L29").Value = 2.52
$ws.Range("M29").Value = 1.06
$ws.Range("N29").Value = 7.3
$ws.Range("O29").Value = 1.3
$ws.Range("P29").Value = 3.2
$ws.Range("Q29").Value = 1.9
$ws.Range("R29").Value = 1.83
$ws.Range("S29").Value = 1.42
$ws.Range("T29").Value = 2.67
$ws.Range("U29").Value = 1.78
$ws.Range("V29").Value = 1.93
$ws.Range("W29").Value = 10.5
$ws.Range("X29").Value = 19.5
$ws.Range("Z29").Value = 55
$ws.Range("AA29").Value = 35
$ws.Range("AB29").Value = 40
$ws.Range("AC29").Value = 7.3
$ws.Range("AD29").Value = 6.7
$ws.Range("AE29").Value = 14.5
$ws.Range("AF29").Value = 70
$ws.Range("AG29").Value = 500
$ws.Range("AH29").Value = 7.3
$ws.Range("AI29").Value = 9.25
$ws.Range("AK29").Value = 16.5
$ws.Range("AL29").Value = 15
$ws.Range("AM29").Value = 26
$ws.Range("AN29").Value = 5.5
$ws.Range("AO29").Value = 21
$ws.Range("AP29").Value = 29
$ws.Range("AQ29").Value = 110
$ws.Range("AR29").Value = 150
$ws.Range("AS29").Value = 400
$ws.Range("AT29").Value = 2.67
$ws.Range("AU29").Value = 7.4
$ws.Range("AV29").Value = 70
$ws.Range("AW29").Value = 3.8
$ws.Range("AX29").Value = 9.75
$ws.Range("AY29").Value = 19
$ws.Range("AZ29").Value = 37
$ws.Range("BA29").Value = 70
$ws.Range("BB29").Value = 250
